$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 7012
$ws.Range("B2").Value = 'Stella Cunha'
$ws.Range("C2").Value = 'Marketing'
$ws.Range("D2").Value = 'Viagem de negócios'
$ws.Range("E2").Value = 5
$ws.Range("G2").Value = 3125.27

$ws.Range("A3").Value = 7998
$ws.Range("B3").Value = 'Giovanna Pereira'
$ws.Range("C3").Value = 'P&D'
$ws.Range("D3").Value = 'Problemas pessoais'
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5631.12

$ws.Range("A4").Value = 32046
$ws.Range("B4").Value = 'Alice da Mota'
$ws.Range("C4").Value = 'Marketing'
$ws.Range("E4").Value = 7
$ws.Range("G4").Value = 11290.13

$ws.Range("A5").Value = 62371
$ws.Range("B5").Value = 'Bruna Nunes'
$ws.Range("C5").Value = 'Recursos Humanos'
$ws.Range("D5").Value = 'Consulta médica'
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 6479.55

$ws.Range("A6").Value = 14173
$ws.Range("B6").Value = 'Kamilly Barros'
$ws.Range("C6").Value = 'Marketing'
$ws.Range("D6").Value = 'Viagem de negócios'
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3119.86

$ws.Range("A7").Value = 23124
$ws.Range("B7").Value = 'Ryan Almeida'
$ws.Range("C7").Value = 'Marketing'
$ws.Range("E7").Value = 6
$ws.Range("G7").Value = 12022.7

$ws.Range("A8").Value = 98115
$ws.Range("B8").Value = 'Ana Beatriz Silva'
$ws.Range("C8").Value = 'Recursos Humanos'
$ws.Range("D8").Value = 'Viagem de negócios'
$ws.Range("E8").Value = 7
$ws.Range("G8").Value = 6885.42

$ws.Range("A9").Value = 95240
$ws.Range("B9").Value = 'Dra. Rafaela Oliveira'
$ws.Range("C9").Value = 'TI'
$ws.Range("D9").Value = 'Consulta médica'
$ws.Range("E9").Value = 5
$ws.Range("G9").Value = 2689.02

$ws.Range("A10").Value = 54398
$ws.Range("B10").Value = 'Miguel da Luz'
$ws.Range("C10").Value = 'Financeiro'
$ws.Range("D10").Value = 'Doença'
$ws.Range("E10").Value = 1
$ws.Range("G10").Value = 4823.06

$ws.Range("A11").Value = 37153
$ws.Range("B11").Value = 'Larissa Lima'
$ws.Range("C11").Value = 'P&D'
$ws.Range("D11").Value = 'Problemas pessoais'
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 8619.24
